$wb = $excel.ActiveWorkbook

# --- Sheet "All Orders": insert a new order row at the top of the data (row 2) ---
$ws = $wb.Worksheets.Item("All Orders")
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = 26
$ws.Range("B2").Value = "2026-01-20 17:36"
$ws.Range("C2").Value = "Priyanka Patil"
$ws.Range("D2").Value = "A-1605"
# Phone number must stay text (it's all digits, so force it with a quote prefix).
$ws.Range("E2").Value = "'9867003224"
$ws.Range("F2").Value = "Appe Chutney x2"
$ws.Range("G2").Value = 120
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"
# Collection date must stay text (pure date string), so force it with a quote prefix.
$ws.Range("J2").Value = "'2026-01-21"
$ws.Range("K2").Value = "09:15"
$ws.Range("L2").Value = "'"
$ws.Range("M2").Value = "'"
$ws.Range("N2").Value = "'"

# --- Sheet "Daily Summary": roll the new order into the 2026-01-20 totals ---
$ds = $wb.Worksheets.Item("Daily Summary")
$ds.Range("B2").Value = 7
$ds.Range("E2").Value = 500
$ds.Range("G2").Value = 330

Write-Output "done"
